$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.708.15"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.99%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.863.18"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.83%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.036"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +2.95%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.36"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.031"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.52%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4405"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.68%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3795"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.62%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07454"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8830"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.77"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.882.31"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -8.10%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.549"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.744"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07206"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.74"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.63%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.037"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009106"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.92%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.032"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.41"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.730.13"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.311"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.44"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.79"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.33%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.944"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.84"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.78%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.997"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.305"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.37"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.29%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09059"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.213"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7687"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.017"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.575"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01984"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05346"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.21%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5192"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.831"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1691"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.840"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.695"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.48%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "109.47"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.85%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.63"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.29%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.726"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4690"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.40%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06426"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.52%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.871"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "39.80"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.43"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.31%  "

$ws.Range("E35").Value = "  +2.72%  "
$ws.Range("E36").Value = "  +3.48%  "
